$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: only the Taxonsorteringsordning (B2) id changes ---
$ws.Range("B2").Value = 78700

# --- Row 3 and Row 4 trade most of their species data, but each also ---
# --- picks up a brand new "B" (Taxonsorteringsordning) value, and the ---
# --- empty "Bestämningsmetod" (AF) marker cell moves from row 3 to row 4. ---

# Row 3 -> becomes what used to be "Lunglav" (Lobaria pulmonaria) entry,
# keeping its own id (A) swapped with row 4's, and a fresh B value.
$ws.Range("A3").Value = 112248267
$ws.Range("B3").Value = 78699
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6458
$ws.Range("F3").Value = "Lunglav"
$ws.Range("G3").Value = "Lobaria pulmonaria"
$ws.Range("H3").Value = "(L.) Hoffm."
$ws.Range("Q3").Value = 522973
$ws.Range("R3").Value = 6909162
# The empty inline-string marker cell in AF3 disappears entirely.
$ws.Range("AF3").ClearContents()

# Row 4 -> becomes what used to be "Källmossa" (Philonotis fontana) entry.
$ws.Range("A4").Value = 112248307
$ws.Range("B4").Value = 92666
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 2412
$ws.Range("F4").Value = "Källmossa"
$ws.Range("G4").Value = "Philonotis fontana"
$ws.Range("H4").Value = "(Hedw.) Brid."
$ws.Range("Q4").Value = 522776
$ws.Range("R4").Value = 6909411
# AF4 gains an (empty) marker cell that previously only existed on row 3.
# (A bare empty-string write clears the cell outright in this engine, so
# nudge an unrelated, default-valued property to force the cell to persist.)
$ws.Range("AF4").Value = ""
$ws.Range("AF4").Font.Bold = $false
